$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value / formula changes ---
# "num blocks per run" (experiment block): 3 -> 4
$ws.Range("B9").Value = 4

# "run length" (experiment block): B10*B13 -> B10*B13+B5
$ws.Range("B14").Formula = "=B10*B13+B5"

# "run length " (localizer block): B24*B26 -> B24*B26+B21
$ws.Range("B27").Formula = "=B24*B26+B21"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 44.7
$ws.Columns.Item(2).ColumnWidth = 9.1

# --- Fonts: Arial -> Calibri (touch only the populated cells, one at a
#     time, so empty cells in the range aren't materialised) ---
$fontCells = @(
  "B1","B2","A3","B3","A4","B4","A5","B5","A8","B8","A9","B9","A10","B10",
  "A13","B13","A14","B14","A18","A19","B19","A20","B20","A21","B21","A23",
  "B23","A24","B24","A26","B26","A27","B27"
)
foreach ($addr in $fontCells) {
  $ws.Range($addr).Font.Name = "Calibri"
}

# --- Selection ---
$ws.Range("B28").Select() | Out-Null
